$d = $word.ActiveDocument

$replacements = @(
    @{old="68÷8=8, 4"; new="46÷2=23, 0"},
    @{old="95÷6=15, 5"; new="60÷8=7, 4"},
    @{old="96÷2=48, 0"; new="81÷7=11, 4"},
    @{old="77÷8=9, 5"; new="71÷7=10, 1"},
    @{old="53÷2=26, 1"; new="48÷5=9, 3"},
    @{old="60÷6=10, 0"; new="82÷5=16, 2"},
    @{old="27÷7=3, 6"; new="88÷8=11, 0"},
    @{old="30÷9=3, 3"; new="63÷9=7, 0"},
    @{old="24÷6=4, 0"; new="95÷8=11, 7"},
    @{old="65÷8=8, 1"; new="27÷5=5, 2"},
    @{old="51÷8=6, 3"; new="19÷9=2, 1"},
    @{old="41÷2=20, 1"; new="11÷2=5, 1"},
    @{old="91÷4=22, 3"; new="67÷4=16, 3"},
    @{old="34÷6=5, 4"; new="27÷8=3, 3"},
    @{old="76÷3=25, 1"; new="12÷4=3, 0"},
    @{old="90÷7=12, 6"; new="22÷7=3, 1"},
    @{old="97÷4=24, 1"; new="48÷7=6, 6"},
    @{old="53÷4=13, 1"; new="16÷9=1, 7"},
    @{old="74÷7=10, 4"; new="88÷7=12, 4"},
    @{old="88÷3=29, 1"; new="63÷5=12, 3"},
    @{old="99÷7=14, 1"; new="81÷2=40, 1"},
    @{old="37÷2=18, 1"; new="25÷9=2, 7"},
    @{old="37÷5=7, 2"; new="99÷2=49, 1"},
    @{old="38÷7=5, 3"; new="37÷7=5, 2"},
    @{old="87÷6=14, 3"; new="99÷3=33, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
